# Adding the changes we made on may 9th
#
# 5 new accelerometer readings were inserted right after the header row
# (becoming new rows 2-6), the previously-existing 20 readings shift down
# to rows 7-26, and 5 more new readings are appended at the end
# (rows 27-31). Net effect: the data block grows from A1:C21 to A1:C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = 30
$arr = New-Object 'object[,]' $rowCount,3

$arr[0,0] = 3.055751519543784; $arr[0,1] = -7.220968450818743; $arr[0,2] = -2.894419597727912
$arr[1,0] = 2.955562557492938; $arr[1,1] = -7.134888444628034; $arr[1,2] = -2.944399050303869
$arr[2,0] = 2.833344757556915; $arr[2,1] = -7.338198423385621; $arr[2,2] = -2.509933024644852
$arr[3,0] = 3.108331612178258; $arr[3,1] = -7.027578847748893; $arr[3,2] = -2.627250722476415
$arr[4,0] = 3.771172387259348; $arr[4,1] = -7.547038768018995; $arr[4,2] = -3.545775217669351
$arr[5,0] = 2.675201347896034; $arr[5,1] = -6.581202302660263; $arr[5,2] = 2.435893160956232
$arr[6,0] = 0.8038081441606844; $arr[6,1] = -7.229396700859067; $arr[6,2] = 1.553160754697669
$arr[7,0] = 1.023633156503948; $arr[7,1] = -4.972137170178557; $arr[7,2] = -3.916851524795791
$arr[8,0] = 0.3411836709295031; $arr[8,1] = -4.24002621429307; $arr[8,2] = -5.904611808913096
$arr[9,0] = -0.6227953433990485; $arr[9,1] = -4.929600954055787; $arr[9,2] = -6.7005569934845
$arr[10,0] = -0.2611346755709001; $arr[10,1] = -5.176580480166844; $arr[10,2] = -11.08990749291009
$arr[11,0] = 7.22266825607841; $arr[11,1] = -7.398067985262175; $arr[11,2] = -7.469464063644459
$arr[12,0] = 17.43506016050075; $arr[12,1] = -11.15922947440832; $arr[12,2] = 8.020842296736639
$arr[13,0] = 4.810774190085277; $arr[13,1] = -6.605375971112934; $arr[13,2] = 0.3875993319920126
$arr[14,0] = 3.291332300220217; $arr[14,1] = -5.354068347385952; $arr[14,2] = 1.089433806283133
$arr[15,0] = 3.138020081179482; $arr[15,1] = -6.129925046648298; $arr[15,2] = 3.955937453678676
$arr[16,0] = 3.258425533771516; $arr[16,1] = -4.712551474571227; $arr[16,2] = 5.692880451679232
$arr[17,0] = -3.766383392470183; $arr[17,1] = -8.431432792118596; $arr[17,2] = 7.765124661581846
$arr[18,0] = 1.177985225405012; $arr[18,1] = -3.925881973334721; $arr[18,2] = -4.60967251232692
$arr[19,0] = 1.052592243467076; $arr[19,1] = -3.166507703917357; $arr[19,2] = -6.040607401302906
$arr[20,0] = -1.074721106461116; $arr[20,1] = -3.308761754206249; $arr[20,2] = -6.224646031856537
$arr[21,0] = -1.158144678388323; $arr[21,1] = -3.9557591165815; $arr[21,2] = -8.491945947919575
$arr[22,0] = 0.7903139420917948; $arr[22,1] = -5.482765521321981; $arr[22,2] = -6.279002168348867
$arr[23,0] = 4.010212659835815; $arr[23,1] = -5.284669637680054; $arr[23,2] = -2.184034883975983
$arr[24,0] = 2.293728096144541; $arr[24,1] = -7.460245260170529; $arr[24,2] = 0.3245020040443973
$arr[25,0] = 1.675990547452653; $arr[25,1] = -6.526311159133912; $arr[25,2] = -0.6771522419793252
$arr[26,0] = 3.116939672401965; $arr[26,1] = -5.273013770580299; $arr[26,2] = -3.711685695818474
$arr[27,0] = 2.304866858891079; $arr[27,1] = -5.722443546567645; $arr[27,2] = -4.230944650513785
$arr[28,0] = 2.380605901990617; $arr[28,1] = -5.519177751881738; $arr[28,2] = -4.144632569381169
$arr[29,0] = 2.426761286599295; $arr[29,1] = -5.405185426984516; $arr[29,2] = -3.328171287264147

$ws.Range("A2:C31").Value = $arr
